$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row for "Joyce Njogu" above the old row 3 (Benson Ambaisi) ---
$ws.Rows("3:3").Insert()

# --- Header row (row 1) ---
$ws.Range("A1").Value = "FULL NAMES"
$ws.Range("B1").Value = "Days Worked"
$ws.Range("C1").Value = "Holiday OT"
$ws.Range("D1").Value = "Workday  O"

# New E1 header cell: copy formatting from D1 (bold/border/centered style) then set its text
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value = "Restday OT"

# --- Row 2: Paul Kamau ---
$ws.Range("B2").Value = 17
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0

# --- Row 3: Joyce Njogu (new) ---
$ws.Range("A3").Value = "Joyce               Njogu               "
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0

# --- Row 4: Benson Ambaisi ---
$ws.Range("B4").Value = 21
$ws.Range("C4").Value = 4
$ws.Range("D4").Value = 76.59999999999999
$ws.Range("E4").Value = 11

# --- Row 5: Ernest Odinga ---
$ws.Range("B5").Value = 5
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0

# --- Row 6: Nicholas Awino ---
$ws.Range("B6").Value = 19
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0

# --- Row 7: Job Mwihia ---
$ws.Range("B7").Value = 21
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0

# --- Row 8: James Ommira ---
$ws.Range("B8").Value = 19
$ws.Range("C8").Value = 4
$ws.Range("D8").Value = 63.1
$ws.Range("E8").Value = 6

# --- Row 9: Virginia Ngure ---
$ws.Range("B9").Value = 23
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0

# --- Row 10: Walter Ojero ---
$ws.Range("B10").Value = 22
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 57.72
$ws.Range("E10").Value = 0

# --- Row 11: Barrack Ogonji ---
$ws.Range("B11").Value = 19
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0
